# Commands To Implement.xlsx - apply status/programmer updates to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Status column (D) updates: mark commands as implemented ("V") ---
# Rows already marked "?" become "V"
$ws.Range("D25").Value = "V"
$ws.Range("D26").Value = "V"
$ws.Range("D56").Value = "V"
$ws.Range("D57").Value = "V"

# Rows that previously had no Status get a "V"
$ws.Range("D27").Value = "V"
$ws.Range("D36").Value = "V"
$ws.Range("D37").Value = "V"
$ws.Range("D38").Value = "V"
$ws.Range("D39").Value = "V"
$ws.Range("D40").Value = "V"
$ws.Range("D41").Value = "V"
$ws.Range("D42").Value = "V"
$ws.Range("D43").Value = "V"
$ws.Range("D44").Value = "V"
$ws.Range("D46").Value = "V"
$ws.Range("D47").Value = "V"
$ws.Range("D48").Value = "V"
$ws.Range("D49").Value = "V"
$ws.Range("D50").Value = "V"
$ws.Range("D51").Value = "V"

# --- Programmer column (E) reassignments to "Maor" ---
$ws.Range("E25").Value = "Maor"
$ws.Range("E26").Value = "Maor"
$ws.Range("E27").Value = "Maor"

# Row 45 gets a lower-case "v" status (new shared string, added after "Maor")
$ws.Range("D45").Value = "v"

# --- Update the view: scroll / select the last touched cell ---
$ws.Range("D57").Select()
